$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.527867078781128
$ws.Range("B1").Value = 1.881496787071228
$ws.Range("C1").Value = 1.630297541618347
$ws.Range("D1").Value = 2.367570638656616
$ws.Range("E1").Value = 3.569114446640015
